# Scheduled runner update: refresh market price / profit columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR item-leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 159.28572
$arr[0,1] = 102.5
$arr[0,2] = 500
$arr[0,3] = 102.5
$arr[0,4] = 500
$arr[0,5] = 66.5
$arr[0,6] = -838
$ws.Range("H9:N9").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 60666.668
$arr[0,1] = 70000
$ws.Range("H21:I21").Value = $arr
$ws.Range("K21").Value = 70000
$ws.Range("M21").Value = -69532

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 60666.668
$arr[0,1] = 70000
$ws.Range("H23:I23").Value = $arr
$ws.Range("K23").Value = 70000
$ws.Range("M23").Value = -69766

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 31253444
$arr[0,1] = 3333.1667
$arr[0,2] = 50003510
$arr[0,3] = 3333.1667
$arr[0,4] = 50003510
$arr[0,5] = -3158.1667
$arr[0,6] = -50003860
$ws.Range("H40:N40").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 8202.799999999999
$arr[0,1] = 9305.25
$ws.Range("H76:I76").Value = $arr
$ws.Range("K76").Value = 9305.25
$ws.Range("M76").Value = -8990.25

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 8202.799999999999
$arr[0,1] = 9305.25
$ws.Range("H79:I79").Value = $arr
$ws.Range("K79").Value = 9305.25
$ws.Range("M79").Value = -8213.25

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 1332.1666
$arr[0,1] = 797.5
$ws.Range("H100:I100").Value = $arr
$ws.Range("K100").Value = 797.5
$ws.Range("M100").Value = -256.5

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 38464520
$arr[0,1] = 71430900
$arr[0,2] = 3750
$arr[0,3] = 214292700
$arr[0,4] = 11250
$arr[0,5] = -214290150
$arr[0,6] = -16350
$ws.Range("H137:N137").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3999.68
$arr[0,1] = 2570.2727
$arr[0,2] = 5122.7856
$arr[0,3] = 7710.8181
$arr[0,4] = 15368.3568
$arr[0,5] = -2570.8181
$arr[0,6] = -25648.3568
$ws.Range("H138:N138").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 3823.1807
$arr[0,1] = 3838.0142
$ws.Range("H32:I32").Value = $arr
$ws.Range("K32").Value = 3838.0142
$ws.Range("M32").Value = -3551.0142

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 13658795
$arr[0,1] = 15913697
$ws.Range("H61:I61").Value = $arr
$ws.Range("K61").Value = 15913697
$ws.Range("M61").Value = -15913485

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3060.2917
$arr[0,1] = 3168.2632
$arr[0,2] = 2650
$arr[0,3] = 3168.2632
$arr[0,4] = 2650
$arr[0,5] = -2294.2632
$arr[0,6] = -4398
$ws.Range("H74:N74").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3060.2917
$arr[0,1] = 3168.2632
$arr[0,2] = 2650
$arr[0,3] = 15841.316
$arr[0,4] = 13250
$arr[0,5] = -11473.316
$arr[0,6] = -21986
$ws.Range("H77:N77").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 1292.1666
$arr[0,1] = 1250.5294
$ws.Range("H102:I102").Value = $arr
$ws.Range("K102").Value = 1250.5294
$ws.Range("M102").Value = 371.4706000000001

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1728309.8
$arr[0,1] = 3475.4888
$arr[0,2] = 7698890
$arr[0,3] = 10426.4664
$arr[0,4] = 23096670
$arr[0,5] = -7896.466400000001
$arr[0,6] = -23101730
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 13658795
$arr[0,1] = 15913697
$ws.Range("H136:I136").Value = $arr
$ws.Range("K136").Value = 47741091
$ws.Range("M136").Value = -47738541

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 924.3333
$arr[0,1] = 1037.3846
$ws.Range("H22:I22").Value = $arr
$ws.Range("K22").Value = 1037.3846
$ws.Range("M22").Value = -864.3846000000001

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 943.8570999999999
$arr[0,1] = 1203.8182
$ws.Range("H80:I80").Value = $arr
$ws.Range("K80").Value = 1203.8182
$ws.Range("M80").Value = -205.8181999999999

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 943.8570999999999
$arr[0,1] = 1203.8182
$ws.Range("H83:I83").Value = $arr
$ws.Range("K83").Value = 6019.090999999999
$ws.Range("M83").Value = -1027.090999999999

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 2409
$arr[0,1] = 2624
$ws.Range("H99:I99").Value = $arr
$ws.Range("K99").Value = 2624
$ws.Range("M99").Value = -1126

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4002355.2
$arr[0,1] = 2352.353
$arr[0,2] = 12502362
$arr[0,3] = 7057.059
$arr[0,4] = 37507086
$arr[0,5] = -4522.059
$arr[0,6] = -37512156
$ws.Range("H134:N134").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 19610914
$arr[0,1] = 27029284
$arr[0,2] = 5223.7856
$arr[0,3] = 27029284
$arr[0,4] = 5223.7856
$arr[0,5] = -27028989
$arr[0,6] = -5813.7856
$ws.Range("H31:N31").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 19610914
$arr[0,1] = 27029284
$arr[0,2] = 5223.7856
$arr[0,3] = 27029284
$arr[0,4] = 5223.7856
$arr[0,5] = -27029082
$arr[0,6] = -5627.7856
$ws.Range("H34:N34").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1982.8918
$arr[0,1] = 1885.1923
$arr[0,2] = 2213.818
$arr[0,3] = 1885.1923
$arr[0,4] = 2213.818
$arr[0,5] = -1682.1923
$arr[0,6] = -2619.818
$ws.Range("H58:N58").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 20142.572
$arr[0,1] = 21199.6
$arr[0,2] = 17500
$arr[0,3] = 21199.6
$arr[0,4] = 17500
$arr[0,5] = -20076.6
$arr[0,6] = -19746
$ws.Range("H86:N86").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 20142.572
$arr[0,1] = 21199.6
$arr[0,2] = 17500
$arr[0,3] = 105998
$arr[0,4] = 87500
$arr[0,5] = -100382
$arr[0,6] = -98732
$ws.Range("H89:N89").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 1552.7858
$arr[0,1] = 1214.8334
$ws.Range("H94:I94").Value = $arr
$ws.Range("K94").Value = 1214.8334
$ws.Range("M94").Value = -763.8334

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 9991.606
$arr[0,1] = 6409.8945
$ws.Range("H99:I99").Value = $arr
$ws.Range("K99").Value = 6409.8945
$ws.Range("M99").Value = -4911.8945

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 9991.606
$arr[0,1] = 6409.8945
$ws.Range("H126:I126").Value = $arr
$ws.Range("K126").Value = 19229.6835
$ws.Range("M126").Value = -16759.6835

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1982.8918
$arr[0,1] = 1885.1923
$arr[0,2] = 2213.818
$arr[0,3] = 5655.5769
$arr[0,4] = 6641.454000000001
$arr[0,5] = -3105.5769
$arr[0,6] = -11741.454
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 18492.834
$ws.Range("J69").Value = 26239.5
$ws.Range("L69").Value = 78718.5
$ws.Range("N69").Value = -80340.5

$ws.Range("H72").Value = 18492.834
$ws.Range("J72").Value = 26239.5
$ws.Range("L72").Value = 236155.5
$ws.Range("N72").Value = -244267.5

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4789191
$arr[0,1] = 2166.5715
$arr[0,2] = 7581622
$arr[0,3] = 6499.7145
$arr[0,4] = 22744866
$arr[0,5] = -4579.7145
$arr[0,6] = -22748706
$ws.Range("H107:N107").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 15007123
$arr[0,1] = 37505004
$ws.Range("H140:I140").Value = $arr
$ws.Range("K140").Value = 112515012
$ws.Range("M140").Value = -112509832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 30002728
$ws.Range("J80").Value = 85717944
$ws.Range("L80").Value = 85717944
$ws.Range("N80").Value = -85719940

$ws.Range("H83").Value = 30002728
$ws.Range("J83").Value = 85717944
$ws.Range("L83").Value = 428589720
$ws.Range("N83").Value = -428599704

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4350196
$arr[0,1] = 2151.3684
$arr[0,2] = 25003406
$arr[0,3] = 6454.1052
$arr[0,4] = 75010218
$arr[0,5] = -3924.1052
$arr[0,6] = -75015278
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 11054
$arr[0,1] = 13332.667
$ws.Range("H22:I22").Value = $arr
$ws.Range("K22").Value = 13332.667
$ws.Range("M22").Value = -13037.667

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 11054
$arr[0,1] = 13332.667
$ws.Range("H27:I27").Value = $arr
$ws.Range("K27").Value = 13332.667
$ws.Range("M27").Value = -13225.667

$ws.Range("H46").Value = 1624.75
$ws.Range("J46").Value = 1599.6666
$ws.Range("L46").Value = 1599.6666
$ws.Range("N46").Value = -1975.6666

$ws.Range("H63").Value = 99624.75
$ws.Range("J63").Value = 99624.75
$ws.Range("L63").Value = 99624.75
$ws.Range("N63").Value = -101122.75

$ws.Range("H66").Value = 99624.75
$ws.Range("J66").Value = 99624.75
$ws.Range("L66").Value = 298874.25
$ws.Range("N66").Value = -306362.25

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 4276284.5
$arr[0,1] = 2813.3
$ws.Range("H93:I93").Value = $arr
$ws.Range("K93").Value = 2813.3
$ws.Range("M93").Value = -1565.3

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 2357.4
$arr[0,1] = 2821.75
$ws.Range("H81:I81").Value = $arr
$ws.Range("K81").Value = 5643.5
$ws.Range("M81").Value = -4582.5

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 2357.4
$arr[0,1] = 2821.75
$ws.Range("H84:I84").Value = $arr
$ws.Range("K84").Value = 28217.5
$ws.Range("M84").Value = -22913.5

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 878.3684
$arr[0,1] = 986.1818
$arr[0,2] = 730.125
$arr[0,3] = 1972.3636
$arr[0,4] = 1460.25
$arr[0,5] = -1431.3636
$arr[0,6] = -2542.25
$ws.Range("H100:N100").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 323571.97
$arr[0,1] = 1015.86957
$ws.Range("H136:I136").Value = $arr
$ws.Range("K136").Value = 3047.60871
$ws.Range("M136").Value = -497.60871
